$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 3366.6191
$ws.Range("I11").Value = 3366.6191
$ws.Range("K11").Value = 3366.6191
$ws.Range("M11").Value = -3226.6191
$ws.Range("H15").Value = 2179473
$ws.Range("I15").Value = 2179473
$ws.Range("K15").Value = 6538419
$ws.Range("M15").Value = -6538250
$ws.Range("H48").Value = 1449.25
$ws.Range("I48").Value = 998.5
$ws.Range("J48").Value = 1900
$ws.Range("K48").Value = 2995.5
$ws.Range("L48").Value = 5700
$ws.Range("M48").Value = -2703.5
$ws.Range("N48").Value = -6284
$ws.Range("H56").Value = 1449.25
$ws.Range("I56").Value = 998.5
$ws.Range("J56").Value = 1900
$ws.Range("K56").Value = 2995.5
$ws.Range("L56").Value = 5700
$ws.Range("M56").Value = -2461.5
$ws.Range("N56").Value = -6768
$ws.Range("H76").Value = 16670916
$ws.Range("I76").Value = 33335498
$ws.Range("J76").Value = 6334.3335
$ws.Range("K76").Value = 33335498
$ws.Range("L76").Value = 6334.3335
$ws.Range("M76").Value = -33335183
$ws.Range("N76").Value = -6964.3335
$ws.Range("H79").Value = 16670916
$ws.Range("I79").Value = 33335498
$ws.Range("J79").Value = 6334.3335
$ws.Range("K79").Value = 33335498
$ws.Range("L79").Value = 6334.3335
$ws.Range("M79").Value = -33334406
$ws.Range("N79").Value = -8518.333500000001
$ws.Range("H133").Value = 119999
$ws.Range("J133").Value = 119999
$ws.Range("L133").Value = 119999
$ws.Range("N133").Value = -130119
$ws.Range("H137").Value = 4390121.5
$ws.Range("I137").Value = 6412423
$ws.Range("J137").Value = 8467.277
$ws.Range("K137").Value = 19237269
$ws.Range("L137").Value = 25401.831
$ws.Range("M137").Value = -19234719
$ws.Range("N137").Value = -30501.831
$ws.Range("H138").Value = 3074.7637
$ws.Range("I138").Value = 1913.16
$ws.Range("J138").Value = 4042.7666
$ws.Range("K138").Value = 5739.48
$ws.Range("L138").Value = 12128.2998
$ws.Range("M138").Value = -599.4800000000005
$ws.Range("N138").Value = -22408.2998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3004724.8
$ws.Range("I32").Value = 3246367.2
$ws.Range("K32").Value = 3246367.2
$ws.Range("M32").Value = -3246080.2
$ws.Range("H39").Value = 6132.25
$ws.Range("I39").Value = 6132.25
$ws.Range("K39").Value = 6132.25
$ws.Range("M39").Value = -5612.25
$ws.Range("H45").Value = 23979.098
$ws.Range("I45").Value = 26313.523
$ws.Range("K45").Value = 26313.523
$ws.Range("M45").Value = -25936.523
$ws.Range("H63").Value = 10034.826
$ws.Range("I63").Value = 2511.6667
$ws.Range("K63").Value = 2511.6667
$ws.Range("M63").Value = -1825.6667
$ws.Range("H66").Value = 10034.826
$ws.Range("I66").Value = 2511.6667
$ws.Range("K66").Value = 12558.3335
$ws.Range("M66").Value = -9126.333500000001
$ws.Range("H88").Value = 1469.5
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 1469.5
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H102").Value = 2361.2
$ws.Range("I102").Value = 2603
$ws.Range("J102").Value = 1998.5
$ws.Range("K102").Value = 2603
$ws.Range("L102").Value = 1998.5
$ws.Range("M102").Value = -981
$ws.Range("N102").Value = -5242.5
$ws.Range("H132").Value = 2568120.2
$ws.Range("I132").Value = 3270654.8
$ws.Range("J132").Value = 8886.643
$ws.Range("K132").Value = 9811964.399999999
$ws.Range("L132").Value = 26659.929
$ws.Range("M132").Value = -9809434.399999999
$ws.Range("N132").Value = -31719.929

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 9694.333000000001
$ws.Range("I54").Value = 2041.5
$ws.Range("K54").Value = 2041.5
$ws.Range("M54").Value = -1557.5
$ws.Range("H86").Value = 101016.45
$ws.Range("I86").Value = 798.9286
$ws.Range("J86").Value = 334857.34
$ws.Range("K86").Value = 798.9286
$ws.Range("L86").Value = 334857.34
$ws.Range("M86").Value = 324.0714
$ws.Range("N86").Value = -337103.34
$ws.Range("H89").Value = 101016.45
$ws.Range("I89").Value = 798.9286
$ws.Range("J89").Value = 334857.34
$ws.Range("K89").Value = 3994.643
$ws.Range("L89").Value = 1674286.7
$ws.Range("M89").Value = 1621.357
$ws.Range("N89").Value = -1685518.7
$ws.Range("H99").Value = 2751.5833
$ws.Range("I99").Value = 2820.3635
$ws.Range("J99").Value = 1995
$ws.Range("K99").Value = 2820.3635
$ws.Range("L99").Value = 1995
$ws.Range("M99").Value = -1322.3635
$ws.Range("N99").Value = -4991
$ws.Range("H105").Value = 62519100
$ws.Range("I105").Value = 66686708
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 66686708
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -66684961
$ws.Range("N105").Value = -8494
$ws.Range("H134").Value = 4299.5117
$ws.Range("I134").Value = 3095.4324
$ws.Range("K134").Value = 9286.297200000001
$ws.Range("M134").Value = -6751.297200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6161.636
$ws.Range("I86").Value = 6444.364
$ws.Range("J86").Value = 5878.909
$ws.Range("K86").Value = 6444.364
$ws.Range("L86").Value = 5878.909
$ws.Range("M86").Value = -5321.364
$ws.Range("N86").Value = -8124.909
$ws.Range("H89").Value = 6161.636
$ws.Range("I89").Value = 6444.364
$ws.Range("J89").Value = 5878.909
$ws.Range("K89").Value = 32221.82
$ws.Range("L89").Value = 29394.545
$ws.Range("M89").Value = -26605.82
$ws.Range("N89").Value = -40626.545
$ws.Range("H107").Value = 2066.25
$ws.Range("I107").Value = 932.8
$ws.Range("J107").Value = 3955.3333
$ws.Range("K107").Value = 932.8
$ws.Range("L107").Value = 3955.3333
$ws.Range("M107").Value = 987.2
$ws.Range("N107").Value = -7795.3333
$ws.Range("H132").Value = 3699.647
$ws.Range("I132").Value = 2617.8096
$ws.Range("J132").Value = 5447.231
$ws.Range("K132").Value = 7853.4288
$ws.Range("L132").Value = 16341.693
$ws.Range("M132").Value = -5323.4288
$ws.Range("N132").Value = -21401.693
$ws.Range("H134").Value = 6796.9165
$ws.Range("I134").Value = 4424.552
$ws.Range("J134").Value = 10417.895
$ws.Range("K134").Value = 13273.656
$ws.Range("L134").Value = 31253.685
$ws.Range("M134").Value = -10738.656
$ws.Range("N134").Value = -36323.685

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6642335
$ws.Range("I4").Value = 7015323
$ws.Range("K4").Value = 21045969
$ws.Range("M4").Value = -21045857
$ws.Range("H14").Value = 2048.818
$ws.Range("I14").Value = 2048.818
$ws.Range("K14").Value = 6146.454000000001
$ws.Range("M14").Value = -5973.454000000001
$ws.Range("H56").Value = 11041.667
$ws.Range("I56").Value = 11041.667
$ws.Range("K56").Value = 11041.667
$ws.Range("M56").Value = -10511.667
$ws.Range("H86").Value = 363
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 363
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H129").Value = 29413410
$ws.Range("I129").Value = 1614
$ws.Range("J129").Value = 45456210
$ws.Range("K129").Value = 4842
$ws.Range("L129").Value = 136368630
$ws.Range("M129").Value = 158
$ws.Range("N129").Value = -136378630
$ws.Range("H132").Value = 2842.5
$ws.Range("I132").Value = 4300
$ws.Range("J132").Value = 2426.0715
$ws.Range("K132").Value = 38700
$ws.Range("L132").Value = 21834.6435
$ws.Range("M132").Value = -36170
$ws.Range("N132").Value = -26894.6435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1675.4546
$ws.Range("I97").Value = 1580.7222
$ws.Range("J97").Value = 2101.75
$ws.Range("K97").Value = 1580.7222
$ws.Range("L97").Value = 2101.75
$ws.Range("M97").Value = -1084.7222
$ws.Range("N97").Value = -3093.75
$ws.Range("H107").Value = 1284.5
$ws.Range("J107").Value = 1472.8125
$ws.Range("L107").Value = 1472.8125
$ws.Range("N107").Value = -5312.8125
$ws.Range("H122").Value = 8050.4443
$ws.Range("I122").Value = 2860.8667
$ws.Range("K122").Value = 8582.6001
$ws.Range("M122").Value = -6132.6001
$ws.Range("H132").Value = 4636.357
$ws.Range("I132").Value = 3060.4707
$ws.Range("K132").Value = 9181.4121
$ws.Range("M132").Value = -6651.4121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 9261438
$ws.Range("I100").Value = 22729004
$ws.Range("J100").Value = 2486.75
$ws.Range("K100").Value = 22729004
$ws.Range("L100").Value = 2486.75
$ws.Range("M100").Value = -22728463
$ws.Range("N100").Value = -3568.75
$ws.Range("H132").Value = 6057.467
$ws.Range("I132").Value = 2712.25
$ws.Range("K132").Value = 8136.75
$ws.Range("M132").Value = -5606.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 24999
$ws.Range("J55").Value = 24999
$ws.Range("L55").Value = 24999
$ws.Range("N55").Value = -25553
$ws.Range("H132").Value = 9094867
$ws.Range("I132").Value = 11906860
$ws.Range("K132").Value = 35720580
$ws.Range("M132").Value = -35718050
$ws.Range("H136").Value = 4083.6667
$ws.Range("I136").Value = 1862
$ws.Range("K136").Value = 5586
$ws.Range("M136").Value = -3036
